# Update symbol list (cryptos.xlsx) - refresh prices, 1h-volume labels, and
# the "Hora" (hour) column from 13 -> 14 for the latest data pull.
#
# NOTE: all affected cells are stored as TEXT in the workbook (inline
# strings), even the numeric-looking Price ("D") and Hora ("G") columns.
# Assigning a plain numeric-looking string to .Value lets Excel's COM layer
# auto-coerce it to a real number, which would change the cell's stored
# type. To keep these as text (matching the original authoring), we set
# .Value2 with a leading apostrophe for anything that looks numeric - the
# same trick Excel itself uses when you type e.g. '14 into a cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - BNB
$ws.Range("D2").Value2 = "'243.77"
$ws.Range("G2").Value2 = "'14"

# Row 3 - OKB
$ws.Range("D3").Value2 = "'24.93"
$ws.Range("G3").Value2 = "'14"

# Row 4 - HuobiToken
$ws.Range("D4").Value2 = "'5.160"
$ws.Range("G4").Value2 = "'14"

# Row 5 - Cronos
$ws.Range("D5").Value2 = "'0.05713"
$ws.Range("G5").Value2 = "'14"

# Row 6 - KuCoinToken
$ws.Range("D6").Value2 = "'6.479"
$ws.Range("G6").Value2 = "'14"

# Row 7 - GateToken
$ws.Range("D7").Value2 = "'2.987"
$ws.Range("G7").Value2 = "'14"

# Row 8 - MXToken
$ws.Range("D8").Value2 = "'0.8103"
$ws.Range("G8").Value2 = "'14"

# Row 9 - FTXToken
$ws.Range("D9").Value2 = "'0.8310"
$ws.Range("G9").Value2 = "'14"

# Row 10 - WazirX
$ws.Range("D10").Value2 = "'0.1333"
$ws.Range("G10").Value2 = "'14"

# Row 11 - MandalaExchangeToken
$ws.Range("D11").Value2 = "'0.06954"
$ws.Range("G11").Value2 = "'14"

# Row 12 - BitrueCoin
$ws.Range("D12").Value2 = "'0.02809"
$ws.Range("G12").Value2 = "'14"

# Row 13 - BitMartToken
$ws.Range("D13").Value2 = "'0.09370"
$ws.Range("G13").Value2 = "'14"

# Row 14 - BitForexToken
$ws.Range("D14").Value2 = "'0.001511"
$ws.Range("G14").Value2 = "'14"

# Row 15 - One
$ws.Range("D15").Value2 = "'0.0005971"
$ws.Range("E15").Value2 = "14OneONE"
$ws.Range("G15").Value2 = "'14"

# Row 16 - TigerCash
$ws.Range("D16").Value2 = "'0.006153"
$ws.Range("G16").Value2 = "'14"

# Row 17 - LEO
$ws.Range("D17").Value2 = "'3.498"
$ws.Range("G17").Value2 = "'14"

# Row 18 - BTSEToken (price unchanged)
$ws.Range("G18").Value2 = "'14"

# Row 19 - BitpandaEcosystemToken
$ws.Range("D19").Value2 = "'0.3195"
$ws.Range("G19").Value2 = "'14"

# Row 20 - LiechtensteinCryptoassetsExchange
$ws.Range("D20").Value2 = "'0.03206"
$ws.Range("G20").Value2 = "'14"

# Row 21 - ProBitToken (price unchanged)
$ws.Range("G21").Value2 = "'14"

# Row 22 - MCDex
$ws.Range("D22").Value2 = "'3.743"
$ws.Range("G22").Value2 = "'14"

# Row 23 - CoinExToken
$ws.Range("D23").Value2 = "'0.04680"
$ws.Range("G23").Value2 = "'14"

# Row 24 - (price unchanged)
$ws.Range("G24").Value2 = "'14"

# Row 25
$ws.Range("D25").Value2 = "'0.001237"
$ws.Range("G25").Value2 = "'14"

# Row 26
$ws.Range("D26").Value2 = "'0.004238"
$ws.Range("G26").Value2 = "'14"

# Row 27
$ws.Range("D27").Value2 = "'0.00009695"
$ws.Range("G27").Value2 = "'14"

# Row 28 - UpBots
$ws.Range("E28").Value2 = "27UpBotsUBXTWorstin24h"
$ws.Range("G28").Value2 = "'14"

# Rows 29-39 - price unchanged, only Hora updates
$ws.Range("G29").Value2 = "'14"
$ws.Range("G30").Value2 = "'14"
$ws.Range("G31").Value2 = "'14"
$ws.Range("G32").Value2 = "'14"
$ws.Range("G33").Value2 = "'14"
$ws.Range("G34").Value2 = "'14"
$ws.Range("G35").Value2 = "'14"
$ws.Range("G36").Value2 = "'14"
$ws.Range("G37").Value2 = "'14"
$ws.Range("G38").Value2 = "'14"
$ws.Range("G39").Value2 = "'14"

# Row 40 - IDEX
$ws.Range("D40").Value2 = "'0.03623"
$ws.Range("G40").Value2 = "'14"

# Row 41 - KickToken
$ws.Range("D41").Value2 = "'0.006267"
$ws.Range("G41").Value2 = "'14"

# Row 42 - BKEXToken
$ws.Range("D42").Value2 = "'0.1050"
$ws.Range("G42").Value2 = "'14"

# Row 43 - CEJI
$ws.Range("D43").Value2 = "'0.002709"
$ws.Range("G43").Value2 = "'14"

# Row 44 - LocalTraders
$ws.Range("D44").Value2 = "'0.007347"
$ws.Range("G44").Value2 = "'14"

# Row 45 - CoinLion
$ws.Range("D45").Value2 = "'0.00005272"
$ws.Range("G45").Value2 = "'14"

# Row 46 - Kangarootoken (price unchanged)
$ws.Range("G46").Value2 = "'14"

# Row 47 - CoinbaseStockToken
$ws.Range("D47").Value2 = "'0.1901"
$ws.Range("G47").Value2 = "'14"

# Rows 48-51 - price unchanged, only Hora updates
$ws.Range("G48").Value2 = "'14"
$ws.Range("G49").Value2 = "'14"
$ws.Range("G50").Value2 = "'14"
$ws.Range("G51").Value2 = "'14"
